$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999991443548086
$ws.Range("E2").Value = 0.9999991443548086

$ws.Range("D3").Value = 0.9355868099533976
$ws.Range("E3").Value = 0.9355868099533976

$ws.Range("D4").Value = [double]"0.9999999999999958"
$ws.Range("E4").Value = [double]"4.218847493575595E-15"

$ws.Range("F4").Value = 5.571281909942627
